$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 125000456
$ws.Range("J5").Value = 200000640
$ws.Range("L5").Value = 200000640
$ws.Range("N5").Value = -200000870

$ws.Range("H33").Value = 336.2857
$ws.Range("I33").Value = 252.5
$ws.Range("K33").Value = 252.5
$ws.Range("M33").Value = -23.5

$ws.Range("H88").Value = 2127.8
$ws.Range("J88").Value = 2272.25
$ws.Range("L88").Value = 2272.25
$ws.Range("N88").Value = -3084.25

$ws.Range("H91").Value = 2127.8
$ws.Range("J91").Value = 2272.25
$ws.Range("L91").Value = 2272.25
$ws.Range("N91").Value = -5080.25

$ws.Range("H125").Value = 2565.7856
$ws.Range("I125").Value = 1345.5
$ws.Range("K125").Value = 12109.5
$ws.Range("M125").Value = -9649.5

$ws.Range("H137").Value = 1397.7222
$ws.Range("I137").Value = 1112.7693
$ws.Range("K137").Value = 3338.3079
$ws.Range("M137").Value = -788.3078999999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5630.5293
$ws.Range("I32").Value = 5333.14
$ws.Range("K32").Value = 5333.14
$ws.Range("M32").Value = -5046.14

$ws.Range("H45").Value = 3060.9644
$ws.Range("I45").Value = 1931.3334
$ws.Range("J45").Value = 3908.1875
$ws.Range("K45").Value = 1931.3334
$ws.Range("L45").Value = 3908.1875
$ws.Range("M45").Value = -1554.3334
$ws.Range("N45").Value = -4662.1875

$ws.Range("H122").Value = 1483.0588
$ws.Range("I122").Value = 1325.9231
$ws.Range("J122").Value = 1993.75
$ws.Range("K122").Value = 3977.7693
$ws.Range("L122").Value = 5981.25
$ws.Range("M122").Value = -1527.7693
$ws.Range("N122").Value = -10881.25

$ws.Range("H132").Value = 4241.0205
$ws.Range("I132").Value = 4030.0908
$ws.Range("K132").Value = 12090.2724
$ws.Range("M132").Value = -9560.2724

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1055.9166
$ws.Range("I20").Value = 999.6667
$ws.Range("J20").Value = 1224.6666
$ws.Range("K20").Value = 999.6667
$ws.Range("L20").Value = 1224.6666
$ws.Range("M20").Value = -752.6667
$ws.Range("N20").Value = -1718.6666

$ws.Range("H86").Value = 6046
$ws.Range("I86").Value = 5920.6
$ws.Range("K86").Value = 5920.6
$ws.Range("M86").Value = -4797.6

$ws.Range("H89").Value = 6046
$ws.Range("I89").Value = 5920.6
$ws.Range("K89").Value = 29603
$ws.Range("M89").Value = -23987

$ws.Range("H134").Value = 8707.102999999999
$ws.Range("I134").Value = 9028.069
$ws.Range("K134").Value = 27084.207
$ws.Range("M134").Value = -24549.207

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 6888.5386
$ws.Range("I86").Value = 6389.75
$ws.Range("J86").Value = 7110.222
$ws.Range("K86").Value = 6389.75
$ws.Range("L86").Value = 7110.222
$ws.Range("M86").Value = -5266.75
$ws.Range("N86").Value = -9356.222

$ws.Range("H89").Value = 6888.5386
$ws.Range("I89").Value = 6389.75
$ws.Range("J89").Value = 7110.222
$ws.Range("K89").Value = 31948.75
$ws.Range("L89").Value = 35551.11
$ws.Range("M89").Value = -26332.75
$ws.Range("N89").Value = -46783.11

$ws.Range("H132").Value = 2715.2368
$ws.Range("I132").Value = 2236.2917
$ws.Range("J132").Value = 3536.2856
$ws.Range("K132").Value = 6708.875100000001
$ws.Range("L132").Value = 10608.8568
$ws.Range("M132").Value = -4178.875100000001
$ws.Range("N132").Value = -15668.8568

$ws.Range("H134").Value = 1610.98
$ws.Range("I134").Value = 1158.7179
$ws.Range("K134").Value = 3476.1537
$ws.Range("M134").Value = -941.1537000000003

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 6175035.5
$ws.Range("I129").Value = 437.2
$ws.Range("K129").Value = 1311.6
$ws.Range("M129").Value = 3688.4

$ws.Range("H131").Value = 12381015
$ws.Range("I131").Value = 6945205
$ws.Range("J131").Value = 16729663
$ws.Range("K131").Value = 20835615
$ws.Range("L131").Value = 50188989
$ws.Range("M131").Value = -20830575
$ws.Range("N131").Value = -50199069

$ws.Range("H137").Value = 61492.21
$ws.Range("I137").Value = 89276.586
$ws.Range("J137").Value = 13861.857
$ws.Range("K137").Value = 267829.758
$ws.Range("L137").Value = 41585.571
$ws.Range("M137").Value = -262729.758
$ws.Range("N137").Value = -51785.571

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4071.2354
$ws.Range("I80").Value = 3567.1
$ws.Range("K80").Value = 3567.1
$ws.Range("M80").Value = -2569.1

$ws.Range("H83").Value = 4071.2354
$ws.Range("I83").Value = 3567.1
$ws.Range("K83").Value = 17835.5
$ws.Range("M83").Value = -12843.5

$ws.Range("H126").Value = 5972
$ws.Range("I126").Value = 6515.6665
$ws.Range("J126").Value = 5156.5
$ws.Range("K126").Value = 19546.9995
$ws.Range("L126").Value = 15469.5
$ws.Range("M126").Value = -17076.9995
$ws.Range("N126").Value = -20409.5

$ws.Range("H132").Value = 2669.1082
$ws.Range("J132").Value = 4835.125
$ws.Range("L132").Value = 14505.375
$ws.Range("N132").Value = -19565.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 614.56757
$ws.Range("I16").Value = 549.89655
$ws.Range("J16").Value = 849
$ws.Range("K16").Value = 549.89655
$ws.Range("L16").Value = 849
$ws.Range("M16").Value = -379.89655
$ws.Range("N16").Value = -1189

$ws.Range("H46").Value = 2115.2
$ws.Range("I46").Value = 1800.5
$ws.Range("K46").Value = 1800.5
$ws.Range("M46").Value = -1612.5

$ws.Range("H55").Value = 1074.5714
$ws.Range("I55").Value = 357.5
$ws.Range("J55").Value = 1515.8462
$ws.Range("K55").Value = 357.5
$ws.Range("L55").Value = 1515.8462
$ws.Range("M55").Value = -184.5
$ws.Range("N55").Value = -1861.8462

$ws.Range("H68").Value = 2037.84
$ws.Range("I68").Value = 2038.6
$ws.Range("K68").Value = 2038.6
$ws.Range("M68").Value = -1289.6

$ws.Range("H70").Value = 28931.334
$ws.Range("I70").Value = 14763
$ws.Range("J70").Value = 43099.668
$ws.Range("K70").Value = 14763
$ws.Range("L70").Value = 43099.668
$ws.Range("M70").Value = -14493
$ws.Range("N70").Value = -43639.668

$ws.Range("H71").Value = 2037.84
$ws.Range("I71").Value = 2038.6
$ws.Range("K71").Value = 10193
$ws.Range("M71").Value = -6449

$ws.Range("H73").Value = 28931.334
$ws.Range("I73").Value = 14763
$ws.Range("J73").Value = 43099.668
$ws.Range("K73").Value = 14763
$ws.Range("L73").Value = 43099.668
$ws.Range("M73").Value = -13827
$ws.Range("N73").Value = -44971.668

$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").Value = ""

$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").Value = ""

$ws.Range("H132").Value = 4013.3713
$ws.Range("I132").Value = 3861.158
$ws.Range("J132").Value = 4194.125
$ws.Range("K132").Value = 11583.474
$ws.Range("L132").Value = 12582.375
$ws.Range("M132").Value = -9053.474
$ws.Range("N132").Value = -17642.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9092
$ws.Range("I62").Value = 9283
$ws.Range("J62").Value = 8948.75
$ws.Range("K62").Value = 9283
$ws.Range("L62").Value = 8948.75
$ws.Range("M62").Value = -8659
$ws.Range("N62").Value = -10196.75

$ws.Range("H65").Value = 9092
$ws.Range("I65").Value = 9283
$ws.Range("J65").Value = 8948.75
$ws.Range("K65").Value = 46415
$ws.Range("L65").Value = 44743.75
$ws.Range("M65").Value = -43295
$ws.Range("N65").Value = -50983.75

$ws.Range("H82").Value = 100000
$ws.Range("J82").Value = 100000
$ws.Range("L82").Value = 100000
$ws.Range("N82").Value = -100766

$ws.Range("H85").Value = 100000
$ws.Range("J85").Value = 100000
$ws.Range("L85").Value = 100000
$ws.Range("N85").Value = -102652

$ws.Range("H132").Value = 2147.3215
$ws.Range("I132").Value = 1936.9131
$ws.Range("J132").Value = 3115.2
$ws.Range("K132").Value = 5810.7393
$ws.Range("L132").Value = 9345.599999999999
$ws.Range("M132").Value = -3280.7393
$ws.Range("N132").Value = -14405.6
